# "Generate Report for Handback"
#
# For each language sheet (zh-cn, de-de) the handback run:
#   - flips the Status column (C) from "Ready for handoff" to
#     "Handed back: in sync with en-US" for every data row,
#   - fills in "Latest Target File" (F) and "Latest Handback File" (G)
#     with the same file (and hyperlink) as "Source File Name" (A) /
#     "Latest Handoff File" (D) respectively, because the handback is in
#     sync with the file that was handed off, and
#   - stamps "Latest Handback DateTime" (H) with the actual handback time
#     (previously the zero/placeholder date 0001-01-01 00:00:00).

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$languages = @(
    @{ Sheet = "zh-cn"; HandbackTime = "2016-03-13 13:01:37" },
    @{ Sheet = "de-de"; HandbackTime = "2016-03-13 13:01:44" }
)

foreach ($lang in $languages) {
    $ws = $wb.Worksheets.Item($lang.Sheet)

    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
    if ($lastRow -lt 2) { $lastRow = 2 }

    for ($row = 2; $row -le $lastRow; $row++) {

        # Source File Name (A) / Latest Handoff File (D) hyperlinks, used
        # to mirror onto the new Latest Target File (F) / Latest Handback
        # File (G) columns.
        $aCell = $ws.Cells.Item($row, 1)
        $dCell = $ws.Cells.Item($row, 4)
        $aLink = $ws.Hyperlinks.Item($aCell.Address())
        $dLink = $ws.Hyperlinks.Item($dCell.Address())

        # Status -> handed back, in sync with en-US
        $ws.Cells.Item($row, 3).Value = $newStatus

        # Latest Target File (F) mirrors Source File Name (A)
        $fCell = $ws.Cells.Item($row, 6)
        $ws.Hyperlinks.Add($fCell, $aLink.Address, [Type]::Missing, [Type]::Missing, $aCell.Value2) | Out-Null

        # Latest Handback File (G) mirrors Latest Handoff File (D)
        $gCell = $ws.Cells.Item($row, 7)
        $ws.Hyperlinks.Add($gCell, $dLink.Address, [Type]::Missing, [Type]::Missing, $dCell.Value2) | Out-Null

        # Latest Handback DateTime (H) gets the real handback timestamp
        $ws.Cells.Item($row, 8).Value = $lang.HandbackTime
    }
}
